# Apply the edit described by the diff:
# 1. Update sigma_010 (sheet2) B2:C12 with refined values
# 2. Update sigma_025 (sheet3) B2:C12 with refined values
# 3. Add a new sheet sigma_050 (sheet4) after sigma_025 with full data

$wb = $excel.ActiveWorkbook

# --- Update existing sheet "sigma_010" ---
$ws010 = $wb.Worksheets.Item("sigma_010")
$ws010.Cells.Item(2, 2).Value = 28.03289724445905
$ws010.Cells.Item(2, 3).Value = 32.03521104997676
$ws010.Cells.Item(3, 2).Value = 28.02402150655753
$ws010.Cells.Item(3, 3).Value = 32.02774222237985
$ws010.Cells.Item(4, 2).Value = 28.06397885516832
$ws010.Cells.Item(4, 3).Value = 32.02566237327959
$ws010.Cells.Item(5, 2).Value = 28.04559960676314
$ws010.Cells.Item(5, 3).Value = 32.03137073415905
$ws010.Cells.Item(6, 2).Value = 28.06544961391636
$ws010.Cells.Item(6, 3).Value = 32.04576543443483
$ws010.Cells.Item(7, 2).Value = 28.01649177841737
$ws010.Cells.Item(7, 3).Value = 32.0604917141047
$ws010.Cells.Item(8, 2).Value = 28.01050955852516
$ws010.Cells.Item(8, 3).Value = 32.0332303979655
$ws010.Cells.Item(9, 2).Value = 28.0358036116714
$ws010.Cells.Item(9, 3).Value = 32.00646902971704
$ws010.Cells.Item(10, 2).Value = 28.02719759364472
$ws010.Cells.Item(10, 3).Value = 32.0295632113321
$ws010.Cells.Item(11, 2).Value = 28.01096102671006
$ws010.Cells.Item(11, 3).Value = 32.02434577756094
$ws010.Cells.Item(12, 2).Value = 28.03329103958331
$ws010.Cells.Item(12, 3).Value = 32.03198519449104

# --- Update existing sheet "sigma_025" ---
$ws025 = $wb.Worksheets.Item("sigma_025")
$ws025.Cells.Item(2, 2).Value = 19.75954491992771
$ws025.Cells.Item(2, 3).Value = 27.84081840500972
$ws025.Cells.Item(3, 2).Value = 19.75833669518123
$ws025.Cells.Item(3, 3).Value = 27.88855561236787
$ws025.Cells.Item(4, 2).Value = 19.7777209926506
$ws025.Cells.Item(4, 3).Value = 27.8581762917387
$ws025.Cells.Item(5, 2).Value = 19.74459938929837
$ws025.Cells.Item(5, 3).Value = 27.8538310412284
$ws025.Cells.Item(6, 2).Value = 19.78019019738306
$ws025.Cells.Item(6, 3).Value = 27.8666740546728
$ws025.Cells.Item(7, 2).Value = 19.76075453026992
$ws025.Cells.Item(7, 3).Value = 27.8591445594299
$ws025.Cells.Item(8, 2).Value = 19.74603659907749
$ws025.Cells.Item(8, 3).Value = 27.86647198654626
$ws025.Cells.Item(9, 2).Value = 19.78123811670433
$ws025.Cells.Item(9, 3).Value = 27.86685564465053
$ws025.Cells.Item(10, 2).Value = 19.77851296963624
$ws025.Cells.Item(10, 3).Value = 27.87266560763447
$ws025.Cells.Item(11, 2).Value = 19.75861105373862
$ws025.Cells.Item(11, 3).Value = 27.88207309219497
$ws025.Cells.Item(12, 2).Value = 19.76455454638676
$ws025.Cells.Item(12, 3).Value = 27.86552662954736

# --- Add new sheet "sigma_050" after "sigma_025" ---
$wsAfter = $wb.Worksheets.Item("sigma_025")
$ws050 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsAfter)
$ws050.Name = "sigma_050"
$ws050.Range("A1").Value = "Rows"
$ws050.Range("B1").Value = "Noisy"
$ws050.Range("C1").Value = "NLM-LBP"
$ws050.Cells.Item(2, 1).Value = 0
$ws050.Cells.Item(2, 2).Value = 14.84047174119649
$ws050.Cells.Item(2, 3).Value = 23.8749351118176
$ws050.Cells.Item(3, 1).Value = 1
$ws050.Cells.Item(3, 2).Value = 14.83749718817578
$ws050.Cells.Item(3, 3).Value = 23.8859531526324
$ws050.Cells.Item(4, 1).Value = 2
$ws050.Cells.Item(4, 2).Value = 14.83905784671976
$ws050.Cells.Item(4, 3).Value = 23.82463860891088
$ws050.Cells.Item(5, 1).Value = 3
$ws050.Cells.Item(5, 2).Value = 14.8387158888531
$ws050.Cells.Item(5, 3).Value = 23.81569066156587
$ws050.Cells.Item(6, 1).Value = 4
$ws050.Cells.Item(6, 2).Value = 14.84219329810063
$ws050.Cells.Item(6, 3).Value = 23.87797679000558
$ws050.Cells.Item(7, 1).Value = 5
$ws050.Cells.Item(7, 2).Value = 14.85473289570397
$ws050.Cells.Item(7, 3).Value = 23.85789115049292
$ws050.Cells.Item(8, 1).Value = 6
$ws050.Cells.Item(8, 2).Value = 14.82586493849476
$ws050.Cells.Item(8, 3).Value = 23.79051817932234
$ws050.Cells.Item(9, 1).Value = 7
$ws050.Cells.Item(9, 2).Value = 14.84753264003469
$ws050.Cells.Item(9, 3).Value = 23.8534617691277
$ws050.Cells.Item(10, 1).Value = 8
$ws050.Cells.Item(10, 2).Value = 14.8408614644039
$ws050.Cells.Item(10, 3).Value = 23.86461532514861
$ws050.Cells.Item(11, 1).Value = 9
$ws050.Cells.Item(11, 2).Value = 14.81055092772705
$ws050.Cells.Item(11, 3).Value = 23.81844855669033
$ws050.Range("A12").Value = "Média"
$ws050.Cells.Item(12, 2).Value = 14.83774788294101
$ws050.Cells.Item(12, 3).Value = 23.84641293057142

Write-Output "done"
